# Auto-generated edit script: updates Leve price/profit columns (H-N)
# across all 8 sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR) to match
# the scheduled price-refresh commit.

$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H100").Value = 15152963
$ws.Range("I100").Value = 41667220
$ws.Range("J100").Value = 1958
$ws.Range("K100").Value = 41667220
$ws.Range("L100").Value = 1958
$ws.Range("M100").Value = -41666679
$ws.Range("N100").Value = -3040
$ws.Range("H133").Value = 15931.429
$ws.Range("J133").Value = 15931.429
$ws.Range("L133").Value = 15931.429
$ws.Range("N133").Value = -26051.429

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2251.2632
$ws.Range("I2").Value = 2144.1538
$ws.Range("J2").Value = 2483.3333
$ws.Range("K2").Value = 2144.1538
$ws.Range("L2").Value = 2483.3333
$ws.Range("M2").Value = -2031.1538
$ws.Range("N2").Value = -2709.3333
$ws.Range("H32").Value = 4620.6
$ws.Range("I32").Value = 3021.8057
$ws.Range("J32").Value = 11015.777
$ws.Range("K32").Value = 3021.8057
$ws.Range("L32").Value = 11015.777
$ws.Range("M32").Value = -2734.8057
$ws.Range("N32").Value = -11589.777
$ws.Range("H45").Value = 1259.5555
$ws.Range("I45").Value = 1012.61536
$ws.Range("K45").Value = 1012.61536
$ws.Range("M45").Value = -635.61536
$ws.Range("H61").Value = 2797.3225
$ws.Range("I61").Value = 1719.1052
$ws.Range("J61").Value = 4504.5
$ws.Range("K61").Value = 1719.1052
$ws.Range("L61").Value = 4504.5
$ws.Range("M61").Value = -1507.1052
$ws.Range("N61").Value = -4928.5
$ws.Range("H116").Value = 2251.2632
$ws.Range("I116").Value = 2144.1538
$ws.Range("J116").Value = 2483.3333
$ws.Range("K116").Value = 2144.1538
$ws.Range("L116").Value = 2483.3333
$ws.Range("M116").Value = 149.8462
$ws.Range("N116").Value = -7071.3333
$ws.Range("H132").Value = 2202.842
$ws.Range("I132").Value = 1813.7812
$ws.Range("J132").Value = 4277.8335
$ws.Range("K132").Value = 5441.3436
$ws.Range("L132").Value = 12833.5005
$ws.Range("M132").Value = -2911.3436
$ws.Range("N132").Value = -17893.5005
$ws.Range("H133").Value = 36666.668
$ws.Range("J133").Value = 36666.668
$ws.Range("L133").Value = 36666.668
$ws.Range("N133").Value = -41726.668
$ws.Range("H136").Value = 2797.3225
$ws.Range("I136").Value = 1719.1052
$ws.Range("J136").Value = 4504.5
$ws.Range("K136").Value = 5157.3156
$ws.Range("L136").Value = 13513.5
$ws.Range("M136").Value = -2607.3156
$ws.Range("N136").Value = -18613.5

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2251.2632
$ws.Range("I3").Value = 2144.1538
$ws.Range("J3").Value = 2483.3333
$ws.Range("K3").Value = 2144.1538
$ws.Range("L3").Value = 2483.3333
$ws.Range("M3").Value = -2030.1538
$ws.Range("N3").Value = -2711.3333
$ws.Range("H20").Value = 1002.3889
$ws.Range("I20").Value = 876.2273
$ws.Range("J20").Value = 1200.6428
$ws.Range("K20").Value = 876.2273
$ws.Range("L20").Value = 1200.6428
$ws.Range("M20").Value = -629.2273
$ws.Range("N20").Value = -1694.6428
$ws.Range("H21").Value = 23000
$ws.Range("J21").Value = 23000
$ws.Range("L21").Value = 23000
$ws.Range("N21").Value = -23472
$ws.Range("H86").Value = 14307
$ws.Range("I86").Value = 1750.25
$ws.Range("J86").Value = 28657.572
$ws.Range("K86").Value = 1750.25
$ws.Range("L86").Value = 28657.572
$ws.Range("M86").Value = -627.25
$ws.Range("N86").Value = -30903.572
$ws.Range("H89").Value = 14307
$ws.Range("I89").Value = 1750.25
$ws.Range("J89").Value = 28657.572
$ws.Range("K89").Value = 8751.25
$ws.Range("L89").Value = 143287.86
$ws.Range("M89").Value = -3135.25
$ws.Range("N89").Value = -154519.86
$ws.Range("H94").Value = 656.5769
$ws.Range("I94").Value = 714.13635
$ws.Range("J94").Value = 340
$ws.Range("K94").Value = 714.13635
$ws.Range("L94").Value = 340
$ws.Range("M94").Value = -263.13635
$ws.Range("N94").Value = -1242
$ws.Range("H99").Value = 1447.7059
$ws.Range("I99").Value = 1363.6364
$ws.Range("J99").Value = 1601.8334
$ws.Range("K99").Value = 1363.6364
$ws.Range("L99").Value = 1601.8334
$ws.Range("M99").Value = 134.3635999999999
$ws.Range("N99").Value = -4597.8334
$ws.Range("H134").Value = 2414.025
$ws.Range("I134").Value = 1770.4117
$ws.Range("K134").Value = 5311.2351
$ws.Range("M134").Value = -2776.2351

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1505.9722
$ws.Range("I58").Value = 981.75
$ws.Range("J58").Value = 3340.75
$ws.Range("K58").Value = 981.75
$ws.Range("L58").Value = 3340.75
$ws.Range("M58").Value = -778.75
$ws.Range("N58").Value = -3746.75
$ws.Range("H108").Value = 38561.332
$ws.Range("J108").Value = 50342
$ws.Range("L108").Value = 50342
$ws.Range("N108").Value = -58022
$ws.Range("H122").Value = 2093.125
$ws.Range("I122").Value = 1249.1666
$ws.Range("K122").Value = 3747.4998
$ws.Range("M122").Value = -1297.4998
$ws.Range("H132").Value = 2584.7585
$ws.Range("I132").Value = 1906.9565
$ws.Range("J132").Value = 5183
$ws.Range("K132").Value = 5720.8695
$ws.Range("L132").Value = 15549
$ws.Range("M132").Value = -3190.8695
$ws.Range("N132").Value = -20609
$ws.Range("H134").Value = 1945
$ws.Range("I134").Value = 1015.36664
$ws.Range("J134").Value = 5431.125
$ws.Range("K134").Value = 3046.09992
$ws.Range("L134").Value = 16293.375
$ws.Range("M134").Value = -511.0999199999997
$ws.Range("N134").Value = -21363.375
$ws.Range("H136").Value = 1505.9722
$ws.Range("I136").Value = 981.75
$ws.Range("J136").Value = 3340.75
$ws.Range("K136").Value = 2945.25
$ws.Range("L136").Value = 10022.25
$ws.Range("M136").Value = -395.25
$ws.Range("N136").Value = -15122.25

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H137").Value = 8421713
$ws.Range("I137").Value = 16670920
$ws.Range("K137").Value = 50012760
$ws.Range("M137").Value = -50007660
$ws.Range("H138").Value = 825.82355
$ws.Range("I138").Value = 825.82355
$ws.Range("J138").Value = 0
$ws.Range("K138").Value = 2477.47065
$ws.Range("L138").Value = 0
$ws.Range("M138").Value = 2662.52935
$ws.Range("N138").ClearContents()

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2600
$ws.Range("I80").Value = 2700
$ws.Range("J80").Value = 2500
$ws.Range("K80").Value = 2700
$ws.Range("L80").Value = 2500
$ws.Range("M80").Value = -1702
$ws.Range("N80").Value = -4496
$ws.Range("H83").Value = 2600
$ws.Range("I83").Value = 2700
$ws.Range("J83").Value = 2500
$ws.Range("K83").Value = 13500
$ws.Range("L83").Value = 12500
$ws.Range("M83").Value = -8508
$ws.Range("N83").Value = -22484
$ws.Range("H102").Value = 2305.3333
$ws.Range("J102").Value = 2221.7144
$ws.Range("L102").Value = 2221.7144
$ws.Range("N102").Value = -5465.7144
$ws.Range("H107").Value = 1085.4286
$ws.Range("I107").Value = 1584.2858
$ws.Range("J107").Value = 586.5714
$ws.Range("K107").Value = 1584.2858
$ws.Range("L107").Value = 586.5714
$ws.Range("M107").Value = 335.7141999999999
$ws.Range("N107").Value = -4426.5714
$ws.Range("H132").Value = 3132.5715
$ws.Range("I132").Value = 2490.9092
$ws.Range("J132").Value = 4218.4614
$ws.Range("K132").Value = 7472.7276
$ws.Range("L132").Value = 12655.3842
$ws.Range("M132").Value = -4942.7276
$ws.Range("N132").Value = -17715.3842
$ws.Range("H137").Value = 34722.5
$ws.Range("J137").Value = 34722.5
$ws.Range("L137").Value = 34722.5
$ws.Range("N137").Value = -44922.5
$ws.Range("H138").Value = 31075
$ws.Range("J138").Value = 31075
$ws.Range("L138").Value = 31075
$ws.Range("N138").Value = -41355
$ws.Range("H139").Value = 29749.5
$ws.Range("J139").Value = 29749.5
$ws.Range("L139").Value = 29749.5
$ws.Range("N139").Value = -40029.5

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 4335.25
$ws.Range("I40").Value = 3166.6667
$ws.Range("J40").Value = 4541.4707
$ws.Range("K40").Value = 3166.6667
$ws.Range("L40").Value = 4541.4707
$ws.Range("M40").Value = -3030.6667
$ws.Range("N40").Value = -4813.4707
$ws.Range("H61").Value = 10934.267
$ws.Range("I61").Value = 11111.556
$ws.Range("J61").Value = 10668.333
$ws.Range("K61").Value = 11111.556
$ws.Range("L61").Value = 10668.333
$ws.Range("M61").Value = -10909.556
$ws.Range("N61").Value = -11072.333
$ws.Range("H82").Value = 1162.625
$ws.Range("I82").Value = 499.66666
$ws.Range("J82").Value = 1560.4
$ws.Range("K82").Value = 499.66666
$ws.Range("L82").Value = 1560.4
$ws.Range("M82").Value = -138.66666
$ws.Range("N82").Value = -2282.4
$ws.Range("H85").Value = 1162.625
$ws.Range("I85").Value = 499.66666
$ws.Range("J85").Value = 1560.4
$ws.Range("K85").Value = 499.66666
$ws.Range("L85").Value = 1560.4
$ws.Range("M85").Value = 748.33334
$ws.Range("N85").Value = -4056.4
$ws.Range("H93").Value = 1688.8
$ws.Range("I93").Value = 1000
$ws.Range("J93").Value = 4444
$ws.Range("K93").Value = 1000
$ws.Range("L93").Value = 4444
$ws.Range("M93").Value = 248
$ws.Range("N93").Value = -6940
$ws.Range("H113").Value = 10934.267
$ws.Range("I113").Value = 11111.556
$ws.Range("J113").Value = 10668.333
$ws.Range("K113").Value = 11111.556
$ws.Range("L113").Value = 10668.333
$ws.Range("M113").Value = -8941.556
$ws.Range("N113").Value = -15008.333

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 69327.734
$ws.Range("I122").Value = 101871.6
$ws.Range("J122").Value = 4240
$ws.Range("K122").Value = 305614.8
$ws.Range("L122").Value = 12720
$ws.Range("M122").Value = -303164.8
$ws.Range("N122").Value = -17620
$ws.Range("H126").Value = 202100.2
$ws.Range("I126").Value = 251875.25
$ws.Range("J126").Value = 3000
$ws.Range("K126").Value = 755625.75
$ws.Range("L126").Value = 9000
$ws.Range("M126").Value = -753155.75
$ws.Range("N126").Value = -13940
$ws.Range("H132").Value = 15629123
$ws.Range("I132").Value = 19234716
$ws.Range("J132").Value = 4890
$ws.Range("K132").Value = 57704148
$ws.Range("L132").Value = 14670
$ws.Range("M132").Value = -19730

Write-Output "Applied scheduled price/profit refresh to all sheets."